$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '72.077.11'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '4.032.31'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '538.81'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = '148.94'
$ws.Range("E6").Value = '  -2.49%  '
$ws.Range("D7").Value = '4.026.14'
$ws.Range("E7").Value = '  -0.29%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  -1.66%  '
$ws.Range("E11").Value = '  -1.25%  '
$ws.Range("D12").Value = '53.16'
$ws.Range("E12").Value = '  +10.48%  '
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("D15").Value = '4.682.28'
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").Value = '4.053.36'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").Value = '14.26'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = '20.63'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D19").Value = '1.20'
$ws.Range("E19").Value = '  -1.14%  '
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("D21").Value = '72.112.52'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").Value = '438.83'
$ws.Range("E22").Value = '  +0.65%  '
$ws.Range("D23").Value = '97.89'
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("D24").Value = '3.50'
$ws.Range("E24").Value = '  -3.77%  '
$ws.Range("D25").Value = '4.28'
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("D26").Value = '14.57'
$ws.Range("E26").Value = '  -1.38%  '
$ws.Range("D27").Value = '4.48'
$ws.Range("E27").Value = '  +28.60%  '
$ws.Range("D28").Value = '11.25'
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("D29").Value = '10.66'
$ws.Range("E29").Value = '  -2.79%  '
$ws.Range("E30").Value = '  +2.12%  '
$ws.Range("D31").Value = '37.11'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").Value = '8.34'
$ws.Range("E32").Value = '  +21.99%  '
$ws.Range("E33").Value = '  +1.42%  '
$ws.Range("D34").Value = '13.53'
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("D35").Value = '49.29'
$ws.Range("E35").Value = '  +14.45%  '
$ws.Range("D36").Value = '680.95'
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("D37").Value = '66.66'
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("E38").Value = '  +3.95%  '
$ws.Range("D39").Value = '0.0₃0909'
$ws.Range("E39").Value = '  +6.96%  '
$ws.Range("B40").Value = 'THORChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D40").Value = '11.34'
$ws.Range("E40").Value = '  +18.29%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.148'
$ws.Range("E41").Value = '  -6.87%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D42").Value = '3.40'
$ws.Range("E42").Value = '  +4.02%  '
$ws.Range("D43").Value = '3.40'
$ws.Range("E43").Value = '  -1.79%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value = '0.0493'
$ws.Range("E46").Value = '  -1.19%  '
$ws.Range("E47").Value = '  -1.90%  '
$ws.Range("D48").Value = '2.63'
$ws.Range("E48").Value = '  -3.43%  '
$ws.Range("E49").Value = '  +1.92%  '
$ws.Range("D50").Value = '3.34'
$ws.Range("E50").Value = '  -3.02%  '
$ws.Range("D51").Value = '0.000280'
$ws.Range("E51").Value = '  +2.20%  '
